# Reference [2] on the "References" slide (last slide, index 30):
#   - The leading "M" of "Marc-Andre" currently shares the same bold/red
#     run as "[2] ". Split it into its own run and restyle it to match
#     the italic/accent1 styling used by the rest of the author name
#     (not bold, italic, accent1 theme color) instead of bold red.
#   - The three trailing runs ". 2021. MPI in " / "Small Bites. " /
#     "PPCES 2021. " (which only differed by a stray missing dirty="0"
#     on the middle one) are merged back into a single run.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(30)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$full = $tr.Text

# --- 1) Split "M" out of the "[2] M" run and restyle it ------------------
$idx = $full.IndexOf("[2] M")
# 1-based character index of the "M" (5th character of "[2] M")
$mChar = $tr.Characters($idx + 5, 1)
$mChar.Font.Color.ObjectThemeColor = 5   # msoThemeColorAccent1 -> accent1
$mChar.Font.Italic = $true
$mChar.Font.Bold = $false

# --- 2) Merge ". 2021. MPI in " + "Small Bites. " + "PPCES 2021. " -------
$full = $tr.Text
$mergedText = ". 2021. MPI in Small Bites. PPCES 2021. "
$idx2 = $full.IndexOf($mergedText)
$mergedRange = $tr.Characters($idx2 + 1, $mergedText.Length)
$mergedRange.Text = $mergedText
